$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.267.72"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.884.03"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "237.79"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "0.4814"
$ws.Range("E7").Value = "  -2.42%  "
$ws.Range("D8").Value = "0.2882"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").Value = "0.06595"
$ws.Range("E9").Value = "  -2.73%  "
$ws.Range("D10").Value = "1.906.34"
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "16.86"
$ws.Range("E11").Value = "  -2.31%  "
$ws.Range("D12").Value = "0.07381"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "5.163"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "87.59"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "0.6601"
$ws.Range("E15").Value = "  -2.35%  "
$ws.Range("D16").Value = "30.260.05"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "13.44"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("D19").Value = "0.000007735"
$ws.Range("E19").Value = "  -3.15%  "
$ws.Range("D20").Value = "5.453"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "2.148.73"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "190.92"
$ws.Range("E23").Value = "  -4.06%  "
$ws.Range("D24").Value = "6.194"
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "9.431"
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("D26").Value = "164.83"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").Value = "18.24"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").Value = "1.935"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").Value = "1.449"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").Value = "4.255"
$ws.Range("E30").Value = "  -2.78%  "
$ws.Range("D31").Value = "'0.09160"
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("D32").Value = "'4.040"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "0.05073"
$ws.Range("E33").Value = "  -3.91%  "
$ws.Range("D34").Value = "0.7357"
$ws.Range("E34").Value = "  -1.35%  "
$ws.Range("D35").Value = "1.148"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("D36").Value = "2.715"
$ws.Range("D37").Value = "0.01825"
$ws.Range("E37").Value = "  -1.62%  "
$ws.Range("D38").Value = "2.645"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("D39").Value = "'0.9180"
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "2.072"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "5.891"
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").Value = "106.06"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "0.4322"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("D45").Value = "0.1366"
$ws.Range("E45").Value = "  -2.84%  "
$ws.Range("D46").Value = "7.632"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").Value = "1.573"
$ws.Range("E47").Value = "  +7.85%  "
$ws.Range("D48").Value = "65.15"
$ws.Range("E48").Value = "  -10.03%  "
$ws.Range("D49").Value = "8.913"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").Value = "34.21"
$ws.Range("E50").Value = "  -3.33%  "
$ws.Range("D51").Value = "0.05764"
$ws.Range("E51").Value = "  -2.25%  "
